$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Date tweaks on a few existing rows (column A)
# ---------------------------------------------------------------------------
$ws.Range("A3").Value  = 45820
$ws.Range("A4").Value  = 45826
$ws.Range("A8").Value  = 45826
$ws.Range("A10").Value = 45820

# ---------------------------------------------------------------------------
# 2) Row 4 - trade closed for a profit
# ---------------------------------------------------------------------------
$ws.Range("G4").Value   = 267.62
$ws.Range("K4").Value   = "Profit"
$ws.Range("L4").Formula = "=(G4-D4)*H4"
$ws.Range("N4").Value   = "Trailing Stop Loss Hit"
$ws.Range("T4").Value   = "Closed"

# ---------------------------------------------------------------------------
# 3) Row 8 - trade closed for a loss (H8 becomes a hard-coded value)
# ---------------------------------------------------------------------------
$ws.Range("G8").Value   = 348.95
$ws.Range("H8").Value   = 176
$ws.Range("K8").Value   = "Loss"
$ws.Range("L8").Formula = "=(G8-D8)*H8"
$ws.Range("N8").Value   = "SL Hit"
$ws.Range("T8").Value   = "Closed"

# ---------------------------------------------------------------------------
# 4) Prep rows 19-24 so they carry the same bordered "blank template" look
#    that row 18 used to have, by copying row 18's original formatting down
#    before row 18 itself gets overwritten with the new trade below.
# ---------------------------------------------------------------------------
$ws.Range("A18:U18").Copy() | Out-Null
$ws.Range("A19:U24").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Give row 18's date cell the same date style used elsewhere in column A.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A18").PasteSpecial(-4122) | Out-Null        # xlPasteFormats
$excel.CutCopyMode = $false

# Row 18 - new trade entry
$ws.Range("A18").Value   = 45814
$ws.Range("B18").Value   = "laurus labs"
$ws.Range("C18").Value   = "Long"
$ws.Range("D18").Value   = 645
$ws.Range("E18").Value   = 579.8
$ws.Range("F18").Value   = 790
$ws.Range("H18").Formula = "=I18/(D18-E18)"
$ws.Range("I18").Value   = 3000
$ws.Range("J18").Formula = "=(F18-D18)/(D18-E18)"

# K18 / L18 / N18 / O18 / P18 end up fully empty (no formatting at all) in
# the final sheet, so clear both contents and formatting on those cells.
$ws.Range("K18").Clear() | Out-Null
$ws.Range("L18").Clear() | Out-Null
$ws.Range("N18").Clear() | Out-Null
$ws.Range("O18").Clear() | Out-Null
$ws.Range("P18").Clear() | Out-Null

$ws.Range("M18").Value = "Daimod patter "
$ws.Range("M18").Style = "Normal"

$ws.Range("Q18").Value = "Weekly Day"
$ws.Range("Q18").Style = "Normal"

$ws.Range("R18").Value = "INR"
$ws.Range("R18").Style = "Normal"

$ws.Range("S18").Formula = "=L18"

$ws.Range("T18").Value = "Active"

$ws.Range("U18").Formula = "=H18*D18"

# ---------------------------------------------------------------------------
# 5) Rows 19-24 stay blank templates (same pattern the old row 18 used to
#    have): H and U reference the (empty) row so they evaluate to #DIV/0!,
#    and S references L (which is empty, so it is 0).
# ---------------------------------------------------------------------------
for ($r = 19; $r -le 24; $r++) {
    $ws.Range("H$r").Formula = "=I$r/(D$r-E$r)"
    $ws.Range("S$r").Formula = "=L$r"
    $ws.Range("U$r").Formula = "=H$r*D$r"
}

$ws.Range("K17").Select() | Out-Null
